$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "cryptos" price/volume snapshot (GitHub Actions refresh).
# Column D ("Price") values that look numeric (e.g. "305.32") must be written
# with a leading single-quote so Excel stores them as TEXT (matching the
# original inline-string cells) instead of auto-converting them to numbers,
# which would lose the exact decimal text (e.g. "24.00" -> 24) and introduce
# floating point artifacts. The Style reset afterwards clears the transient
# quote-prefix formatting so the cell's style index is left untouched.

$ws.Range("D2").Value = '''41.850.41'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '''2.272.06'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''305.32'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").Value = '''93.08'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").Value = '''6.68'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").Value = '''2.623.40'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '''14.34'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").Value = '''2.275.48'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '''0.782'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("D18").Value = '''41.777.98'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = '''12.78'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +4.55%  '
$ws.Range("D20").Value = '''0.0₃0916'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").Value = '''5.97'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = '''68.02'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").Value = '''244.03'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''24.00'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("E29").Value = '  -8.88%  '
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").Value = '''159.17'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").Value = '''5.36'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +3.63%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").Value = '''17.19'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '''2.006.30'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -2.33%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''19.60'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("D45").Value = '''0.0282'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").Value = '''10.27'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +1.52%  '
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").Value = '''53.47'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").Value = '''72.67'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("D51").Value = '''1.15'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.25%  '
